$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) Merge the 4 separate runs of the "10/12/24" heading into a
#    single run reading "10/12/24" (keeps the first run - and its
#    empty rPr - and removes the trailing runs' duplicate text).
# -----------------------------------------------------------------
$datePar = $d.Paragraphs(16)
$pStart = $datePar.Range.Start
$pEnd = $datePar.Range.End - 1   # exclude the paragraph mark

$firstRun = $d.Range($pStart, $pStart + 2)   # "10"
if ($firstRun.Text -ne "10") {
    throw "Unexpected heading text: [$($firstRun.Text)]"
}

$restRange = $d.Range($pStart + 2, $pEnd)    # "/12/24"
$restRange.Delete()

$insertPoint = $d.Range($pStart + 2, $pStart + 2)
$insertPoint.InsertAfter("/12/24")

# -----------------------------------------------------------------
# 2) After the "Risolto finalmente..." paragraph, add:
#      - a blank Normal paragraph (spacing before=0 / after=160)
#      - a new Heading2 paragraph "29/1/25"
#      - a new Normal paragraph with the K-means note
# -----------------------------------------------------------------
$solvedPar = $d.Paragraphs(17)
if ($solvedPar.Range.Text -notmatch "Risolto finalmente") {
    throw "Unexpected paragraph 17 text: [$($solvedPar.Range.Text)]"
}

# --- blank paragraph -------------------------------------------------
$solvedPar.Range.InsertParagraphAfter() | Out-Null
$d = $word.ActiveDocument
$blankPar = $d.Paragraphs(18)
$blankPar.Style = "Normal"
$blankPar.SpaceBefore = 0
$blankPar.SpaceAfter = 8
$blankPar.Range.Font.Italic = 0
$blankPar.Range.Font.ItalicBi = 0

# --- "29/1/25" heading -------------------------------------------------
$blankPar.Range.InsertParagraphAfter() | Out-Null
$d = $word.ActiveDocument
$newHeading = $d.Paragraphs(19)
$newHeading.Style = "Heading2"
$newHeading.Range.Font.Italic = 0
$newHeading.Range.Font.ItalicBi = 0
$hStart = $newHeading.Range.Start
$collapsed = $d.Range($hStart, $hStart)
$collapsed.InsertAfter("29")
$d = $word.ActiveDocument
$newHeading = $d.Paragraphs(19)
$afterFirst = $newHeading.Range.Start + 2
$collapsed2 = $d.Range($afterFirst, $afterFirst)
$collapsed2.InsertAfter("/1/2")
$d = $word.ActiveDocument
$newHeading = $d.Paragraphs(19)
$afterSecond = $newHeading.Range.Start + 2 + 4
$collapsed3 = $d.Range($afterSecond, $afterSecond)
$collapsed3.InsertAfter("5")

# --- K-means paragraph -------------------------------------------------
$d = $word.ActiveDocument
$newHeading = $d.Paragraphs(19)
$newHeading.Range.InsertParagraphAfter() | Out-Null
$d = $word.ActiveDocument
$newNormal = $d.Paragraphs(20)
$newNormal.Style = "Normal"
$newNormal.Range.Font.Italic = 0
$newNormal.Range.Font.ItalicBi = 0
$nStart = $newNormal.Range.Start
$nCollapsed = $d.Range($nStart, $nStart)
$nCollapsed.InsertAfter("Aggiunta codice per creazione del dataset di training per algoritmo K-means. Il prossimo step è quello di creare il dataset e scrivere il codice su Colab per il modello di classificazione.")

Write-Output "done"
